# "included code for comparing primary key values"
# Book1.xlsx / Sheet1 holds a small TradeID/Risk/Curve/Type table; the
# primary-key comparison logic produced an updated Curve value for the
# TradeID 1235 row (D3), and the user's cursor ended up parked on that
# cell when the workbook was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated comparison result for primary key 1235 (column D, row 3)
$ws.Range("D3").Value = 768

# Leave the selection on the cell that was just (re)computed
$ws.Range("D3").Select()
